# Auto-generated Excel COM-interop script to apply cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.130.61'
$ws.Range('E2').Value = '  +3.03%  '
$ws.Range('D3').Value = '3.462.96'
$ws.Range('E3').Value = '  +2.35%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''585.02'
$ws.Range('E5').Value = '  +5.20%  '
$ws.Range('D6').Value = '''191.04'
$ws.Range('E6').Value = '  +9.13%  '
$ws.Range('D7').Value = '''0.633'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '3.455.38'
$ws.Range('E8').Value = '  +2.31%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E10').Value = '  +0.17%  '
$ws.Range('D11').Value = '''0.649'
$ws.Range('E11').Value = '  +1.84%  '
$ws.Range('D12').Value = '''57.56'
$ws.Range('E12').Value = '  +7.15%  '
$ws.Range('E13').Value = '  -0.55%  '
$ws.Range('D14').Value = '''9.54'
$ws.Range('E14').Value = '  +3.72%  '
$ws.Range('D15').Value = '4.013.00'
$ws.Range('E15').Value = '  +2.56%  '
$ws.Range('D16').Value = '''18.98'
$ws.Range('E16').Value = '  +3.68%  '
$ws.Range('D17').Value = '3.464.91'
$ws.Range('E17').Value = '  +2.37%  '
$ws.Range('D18').Value = '67.139.45'
$ws.Range('E18').Value = '  +3.22%  '
$ws.Range('D19').Value = '''12.17'
$ws.Range('E19').Value = '  +2.22%  '
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('E21').Value = '  +2.98%  '
$ws.Range('D22').Value = '''481.90'
$ws.Range('E22').Value = '  +6.06%  '
$ws.Range('E23').Value = '  +8.56%  '
$ws.Range('D24').Value = '''16.82'
$ws.Range('E24').Value = '  +18.33%  '
$ws.Range('E25').Value = '  +7.40%  '
$ws.Range('D26').Value = '''90.31'
$ws.Range('E26').Value = '  +3.17%  '
$ws.Range('E27').Value = '  +4.40%  '
$ws.Range('E28').Value = '  +2.79%  '
$ws.Range('D29').Value = '''9.08'
$ws.Range('E29').Value = '  +4.14%  '
$ws.Range('D30').Value = '''31.44'
$ws.Range('E30').Value = '  +0.56%  '
$ws.Range('D31').Value = '''7.49'
$ws.Range('E31').Value = '  +14.63%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = '''11.86'
$ws.Range('E32').Value = '  +3.43%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = '''603.38'
$ws.Range('E33').Value = '  +4.43%  '
$ws.Range('D34').Value = '''64.44'
$ws.Range('E34').Value = '  +2.12%  '
$ws.Range('E35').Value = '  +4.42%  '
$ws.Range('D36').Value = '''0.149'
$ws.Range('E36').Value = '  +5.00%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('E38').Value = '  +4.95%  '
$ws.Range('E39').Value = '  +5.15%  '
$ws.Range('E40').Value = '  -4.97%  '
$ws.Range('D41').Value = '0.0₃0759'
$ws.Range('E41').Value = '  +2.32%  '
$ws.Range('D42').Value = '3.228.94'
$ws.Range('E42').Value = '  +4.21%  '
$ws.Range('D43').Value = '''2.96'
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '''0.0434'
$ws.Range('E44').Value = '  +3.88%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '''2.90'
$ws.Range('E45').Value = '  +28.53%  '
$ws.Range('D46').Value = '''2.58'
$ws.Range('E46').Value = '  +4.80%  '
$ws.Range('D47').Value = '''3.22'
$ws.Range('E47').Value = '  +1.28%  '
$ws.Range('E48').Value = '  +1.10%  '
$ws.Range('D49').Value = '''8.75'
$ws.Range('E49').Value = '  +5.42%  '
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('D51').Value = '''3.22'
$ws.Range('E51').Value = '  +7.62%  '
